$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the merge: clear the values that had been (re)added in column D.
# D3, D5, D7 keep their existing style (s="1"); just drop the values.
$ws.Range("D3").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D7").ClearContents()

# D2 also loses its date-style formatting (s="4" -> s="3"), matching the
# plain, non-date style already used by its neighbour E2.
$ws.Range("D2").ClearContents()
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active cell selection from E9 to E11.
$ws.Range("E11").Select()
